# SolarApiLocations.xlsx edit
# Adds NumPanels / YearlyEnergy / SolarArea data (columns J, K, L) for rows 63-182,
# and disambiguates the Name column for six duplicate-named restaurant rows
# (Slutty Vegan, Sublime Doughnuts, Toast on/On Lenox) by appending their
# distinguishing location to the name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Disambiguate duplicate restaurant names (column A) -------------------
# Order matters: setting them in this sequence reproduces the shared-string
# append order used by the author (McDonough, Edgewood, 10th Street,
# Briarcliff, 14th Street, Lenox).
$ws.Range("A130").Value = "Slutty Vegan McDonough"
$ws.Range("A129").Value = "Slutty Vegan Edgewood"

$ws.Range("A140").Value = "Sublime Doughnuts 10th Street"
$ws.Range("A141").Value = "Sublime Doughnuts Briarcliff"

$ws.Range("A167").Value = "Toast on Lenox 14th Street"
$ws.Range("A168").Value = "Toast On Lenox Lenox"

# --- Fill in NumPanels (J), YearlyEnergy (K), SolarArea (L) ---------------
$data = @(
    @(63, 52, 18737.393, 102.10485),
    @(64, 80, 49429.008000000002, 157.08438000000001),
    @(65, 114, 69919.48, 223.84526),
    @(66, 81, 48009.402000000002, 159.04794000000001),
    @(67, 110, 54934.074000000001, 215.99102999999999),
    @(68, 71, 33310.959999999999, 139.41239999999999),
    @(69, 142, 82909.125, 278.82479999999998),
    @(70, 502, 316975.15999999997, 985.70450000000005),
    @(71, 597, 344774.16, 1172.2422999999999),
    @(72, 97, 49309.934000000001, 190.46483000000001),
    @(73, 310, 190421.75, 608.702),
    @(74, 26, 14267.599, 51.052424999999999),
    @(75, 34, 20253.280999999999, 66.760863999999998),
    @(76, 187, 104780.26, 367.18475000000001),
    @(77, 118, 60355.207000000002, 231.69947999999999),
    @(78, 53, 32035.115000000002, 104.068405),
    @(79, 275, 166204.03, 539.97760000000005),
    @(80, 42, 25216.03, 82.469309999999993),
    @(81, 156, 86055.78, 306.31454000000002),
    @(82, 177, 103876.41, 347.54921999999999),
    @(83, 47, 26713.99, 92.287080000000003),
    @(84, 1614, 899587.94, 3169.1774999999998),
    @(85, 31, 14615.36, 60.870199999999997),
    @(86, 54, 32381.828000000001, 106.03196),
    @(87, 31, 14615.36, 60.870199999999997),
    @(88, 701, 379896.97, 1376.4519),
    @(89, 163, 96094.69, 320.05945000000003),
    @(90, 226, 131312.47, 443.76339999999999),
    @(91, 34, 12983.191999999999, 66.760863999999998),
    @(92, 86, 47995.023000000001, 168.86572000000001),
    @(93, 144, 86613.85, 282.75189999999998),
    @(94, 242, 149661, 475.18027000000001),
    @(95, 45, 22448.728999999999, 88.359970000000004),
    @(96, 131, 74785.554999999993, 257.22568000000001),
    @(97, 444, 217293, 871.81835999999998),
    @(98, 69, 32885.233999999997, 135.48528999999999),
    @(99, 159, 97232.09, 312.20522999999997),
    @(100, 51, 27218.678, 100.1413),
    @(101, 62, 34683.546999999999, 121.74039999999999),
    @(102, 83, 50394.836000000003, 162.97505000000001),
    @(103, 57, 30222.206999999999, 111.92263),
    @(104, 91, 51122.745999999999, 178.68349000000001),
    @(105, 644, 365734.03, 1264.5292999999999),
    @(106, 405, 206903.03, 795.23974999999996),
    @(107, 90, 54021.203000000001, 176.71994000000001),
    @(108, 39, 21568.734, 76.578636000000003),
    @(109, 77, 43016.51, 151.19372999999999),
    @(110, 302, 185321.34, 592.99360000000001),
    @(111, 178, 109520.19, 349.51276000000001),
    @(112, 59, 27434.502, 115.84974),
    @(113, 97, 57024.008000000002, 190.46483000000001),
    @(114, 276, 172367.83, 541.94115999999997),
    @(115, 85, 51610.805, 166.90216000000001),
    @(116, 21, 10737.646000000001, 41.234653000000002),
    @(117, 27, 10997.812, 53.015979999999999),
    @(118, 105, 48861.226999999999, 206.17326),
    @(119, 44, 18556.822, 86.396416000000002),
    @(120, 59, 33454.495999999999, 115.84974),
    @(121, 88, 48923.09, 172.79283000000001),
    @(122, 73, 43358.315999999999, 143.33950999999999),
    @(123, 143, 71239.039999999994, 280.78832999999997),
    @(124, 104, 61771.453000000001, 204.2097),
    @(125, 42, 24369.455000000002, 82.469309999999993),
    @(126, 40, 18206.252, 78.542190000000005),
    @(127, 57, 30222.206999999999, 111.92263),
    @(128, 75, 41834.188000000002, 147.26661999999999),
    @(129, 138, 81442.304999999993, 270.97057999999998),
    @(130, 348, 196594.86, 683.31709999999998),
    @(131, 66, 40822.065999999999, 129.59461999999999),
    @(132, 42, 21377.643, 82.469309999999993),
    @(133, 136, 81245.7, 267.04345999999998),
    @(134, 66, 41587.042999999998, 129.59461999999999),
    @(135, 732, 446733.22, 1437.3221000000001),
    @(136, 1027, 628114.80000000005, 2016.5708),
    @(137, 48, 28080.412, 94.250630000000001),
    @(138, 24, 13332.897999999999, 47.125317000000003),
    @(139, 147, 84517.41, 288.64258000000001),
    @(140, 584, 330379.78000000003, 1146.7161000000001),
    @(141, 597, 344514.25, 1172.2422999999999),
    @(142, 50, 31221.759999999998, 98.17774),
    @(143, 540, 333510.13, 1060.3196),
    @(144, 39, 17985.006000000001, 76.578636000000003),
    @(145, 97, 57024.008000000002, 190.46483000000001),
    @(146, 549, 324342.34000000003, 1077.9916000000001),
    @(147, 7, 4212.2744000000002, 13.7448845),
    @(148, 487, 296738.40000000002, 956.25120000000004),
    @(149, 61, 36085.055, 119.77685),
    @(150, 30, 16912.851999999999, 58.906647),
    @(151, 24, 13748.722, 47.125317000000003),
    @(152, 76, 43858.065999999999, 149.23016000000001),
    @(153, 94, 55848.563000000002, 184.57416000000001),
    @(154, 325, 155863.78, 638.15533000000005),
    @(155, 135, 74776.679999999993, 265.07990000000001),
    @(156, 201, 96264.54, 394.67453),
    @(157, 132, 76885.64, 259.18923999999998),
    @(158, 1129, 661070.80000000005, 2216.8535000000002),
    @(159, 1384, 729206.2, 2717.5598),
    @(160, 354, 208539.61, 695.09844999999996),
    @(161, 530, 323442.2, 1040.6840999999999),
    @(162, 412, 253523.45, 808.9846),
    @(163, 69, 41765.555, 135.48528999999999),
    @(164, 12, 3325.96, 23.562657999999999),
    @(165, 85, 53540.62, 166.90216000000001),
    @(166, 167, 94804.22, 327.91367000000002),
    @(167, 5, 2267.4456, 9.817774),
    @(168, 55, 24537.603999999999, 107.995514),
    @(169, 1740, 1078862.8, 3416.5853999999999),
    @(170, 180, 97480.733999999997, 353.43988000000002),
    @(171, 865, 495516.88, 1698.4749999999999),
    @(172, 301, 183728.2, 591.03),
    @(173, 32, 18315.053, 62.833754999999996),
    @(174, 50, 23862.947, 98.17774),
    @(175, 321, 179966.1, 630.30110000000002),
    @(176, 59, 33332.074000000001, 115.84974),
    @(177, 360, 215623.97, 706.87976000000003),
    @(178, 65, 23114.36, 127.63106500000001),
    @(179, 47, 27945.888999999999, 92.287080000000003),
    @(180, 47, 27682.567999999999, 92.287080000000003),
    @(181, 58, 27227.153999999999, 113.886185),
    @(182, 159, 90219.6, 312.20522999999997)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 10).Value = $row[1]
    $ws.Cells.Item($r, 11).Value = $row[2]
    $ws.Cells.Item($r, 12).Value = $row[3]
}
